$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.062.08'
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").Value = '3.313.44'
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Formula = '="588.98"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").Formula = '="185.11"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +2.12%  '
$ws.Range("D7").Formula = '="0.999"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +2.35%  '
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("D10").Formula = '="6.70"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("D11").Formula = '="0.425"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +2.58%  '
$ws.Range("D12").Value = '3.885.58'
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").Formula = '="0.138"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Formula = '="29.39"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +4.98%  '
$ws.Range("D15").Value = '68.997.81'
$ws.Range("E15").Value = '  +2.49%  '
$ws.Range("E16").Value = '  +3.40%  '
$ws.Range("D17").Value = '3.343.20'
$ws.Range("E17").Value = '  +3.37%  '
$ws.Range("D18").Formula = '="5.88"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").Formula = '="13.71"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +2.12%  '
$ws.Range("D20").Formula = '="393.15"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +4.80%  '
$ws.Range("D21").Formula = '="7.80"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("D22").Formula = '="71.96"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("E24").Value = '  +3.18%  '
$ws.Range("E25").Value = '  +2.07%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Formula = '="9.82"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Formula = '="0.189"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +4.80%  '
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").Formula = '="5.82"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +2.61%  '
$ws.Range("D30").Formula = '="2.00"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("D31").Formula = '="23.21"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("E32").Value = '  +4.54%  '
$ws.Range("D33").Formula = '="7.21"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +5.17%  '
$ws.Range("D34").Formula = '="0.999"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  +4.14%  '
$ws.Range("D36").Formula = '="163.37"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("E37").Value = '  +3.17%  '
$ws.Range("D38").Formula = '="0.841"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D39").Formula = '="26.67"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").Formula = '="4.62"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +4.20%  '
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("D42").Formula = '="6.63"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -2.39%  '
$ws.Range("D43").Formula = '="41.67"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +2.72%  '
$ws.Range("D44").Formula = '="0.0696"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +3.58%  '
$ws.Range("D45").Formula = '="25.48"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").Formula = '="345.71"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -4.88%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.637.67'
$ws.Range("E47").Value = '  -2.63%  '
$ws.Range("D48").Formula = '="0.0286"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +2.77%  '
$ws.Range("D49").Formula = '="32.39"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +5.42%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Formula = '="6.33"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +3.25%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Formula = '="1.00"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +0.64%  '
